$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the data set. It is inserted as row 24,
# pushing every existing row from 24 downward down by one (row 75 -> 76).
$ws.Rows(24).Insert()

$ws.Range("A24").Value = 9
$ws.Range("B24").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 44799
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = 100112029
$ws.Range("G24").Value = "Orégano"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 16
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 20000
$ws.Range("N24").Value = "$/docena de atados"
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 6667
$ws.Range("Q24").Value = 3
$ws.Range("R24").Value = "Hortaliza"
